# Coverage-Payer.xlsx — "Mapping en commentaar Astrid"
# Re-maps column B (the Payer-side mapping) against the fixed column A
# (EHDSCoverage fields), adds a vertical-top alignment to the whole table,
# widens column B slightly, and drops the two now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlTop
$xlTop = -4160

# --- 1. Column B content: re-point each EHDSCoverage row at the right
#        Payer element (some rows now intentionally blank, "Payer" is
#        reused for both EHDSCoverage.payor and EHDSCoverage.identifier-ish
#        rows, and the InsuranceCompany sub-rows are reordered). ---
$ws.Range("B2").Value  = "Payer"
$ws.Range("B3").Value  = ""
$ws.Range("B4").Value  = "Payer.InsuranceCompany.Insurance.InsuranceType"
$ws.Range("B5").Value  = ""
$ws.Range("B6").Value  = "Payer"
$ws.Range("B7").Value  = ""
$ws.Range("B8").Value  = "Payer.InsuranceCompany.InsurantNumber"
$ws.Range("B9").Value  = ""
$ws.Range("B10").Value = "Payer.PayerPerson"
$ws.Range("B11").Value = "Payer.PayerPerson.PayerName"
$ws.Range("B12").Value = "Payer.PayerPerson.BankInformation"
$ws.Range("B13").Value = "Payer.PayerPerson.BankInformation.BankName"
$ws.Range("B14").Value = "Payer.PayerPerson.BankInformation.BankCode"
$ws.Range("B15").Value = "Payer.PayerPerson.BankInformation.AccountNumber"
$ws.Range("B16").Value = "Payer.InsuranceCompany"
$ws.Range("B17").Value = "Payer.InsuranceCompany.Insurance"
$ws.Range("B18").Value = "Payer.InsuranceCompany.OrganizationName"
$ws.Range("B19").Value = "Payer.InsuranceCompany.IdentificationNumber"
$ws.Range("B20").Value = "Payer.InsuranceCompany.Insurance.StartDateTime"
$ws.Range("B21").Value = "Payer.InsuranceCompany.Insurance.EndDateTime"
$ws.Range("B22").Value = "Payer.AddressInformation"
$ws.Range("B23").Value = "Payer.ContactInformation"

# B24 used to hold "Payer.ContactInformation" (now moved up to B23) — clear it.
$ws.Range("B24").ClearContents()

# --- 2. Alignment: vertical-top across the whole mapping block.
#        Row 9 used to be the first "indented" (Calibri) row; it's now a
#        blank separator that matches the plain header rows above it, so
#        drop its Calibri override back to the sheet's default font first.
#        Then do the plain/default-font cells so the new "default font +
#        top" style lands at cellXfs index 3 (matching the target file),
#        and finally the still-Calibri "indent" column-A cells so their
#        "Calibri + top" style lands at index 4. ---
$ws.Range("A9").ClearFormats()
$ws.Range("A2:B9").VerticalAlignment = $xlTop
$ws.Range("B10:B23").VerticalAlignment = $xlTop
$ws.Range("A10:A23").VerticalAlignment = $xlTop

# --- 3. Column B a hair wider, dimension shrink (drop now-empty trailing
#        rows 39:40), and leave the selection where Astrid's edit did. ---
$ws.Range("B1").ColumnWidth = 53.33
$ws.Rows("39:40").Delete()
$ws.Range("B24").Select() | Out-Null
